$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B12 (Enabled column for transport:OpenStreetMap) from "no" to "yes"
$ws.Range("B12").Value = "yes"

# Update selected cell in the sheet view
$ws.Range("C22").Select()
